# ---------------------------------------------------------------------------
# "Add files via upload" - re-creates the edits from the commit:
#   * Sheet1 -> "Tables", Sheet2 -> "Process Flow"
#   * New A16:C23 relationship mini-table (Customer/Account/Address/Request/
#     NetBanking/ATM/Admin <-> cardinality), styled like the other table
#     headers (bold + thin box border, i.e. the format already used by B6)
#   * "Net Banking" table renamed to "NetBanking" (N19)
#   * Email Id gets a new "CHECK(valid)" constraint (L13)
#   * Branch Code / INT / FOREIGN KEY (F22:H22) re-coloured red
#   * Columns B & C widened to fit the new "NetBanking" / "One to Many" text
#   * Selection moved to A27
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsTables = $wb.Worksheets.Item(1)
$wsFlow   = $wb.Worksheets.Item(2)

# 1) Rename the worksheets -------------------------------------------------
$wsTables.Name = "Tables"
$wsFlow.Name   = "Process Flow"

# 2) New relationship table in A16:C23 -------------------------------------
#    Copy the existing bold/bordered header format (used by B6) onto the
#    whole block first, then fill in the text - this re-uses the current
#    style (fontId=1 bold, borderId=1 thin box) instead of minting a new one.
$wsTables.Range("B6").Copy()
$wsTables.Range("A16:C23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsTables.Range("A16").Value = "Customer"
$wsTables.Range("B16").Value = "Account"
$wsTables.Range("C16").Value = "One to Many"

$wsTables.Range("B17").Value = "Address"
$wsTables.Range("C17").Value = "One to Many"

$wsTables.Range("B18").Value = "Request"
$wsTables.Range("C18").Value = "One to One"

$wsTables.Range("B19").Value = "NetBanking"
$wsTables.Range("C19").Value = "One to One"

$wsTables.Range("A21").Value = "Account"
$wsTables.Range("B21").Value = "ATM"
$wsTables.Range("C21").Value = "One to Many"

$wsTables.Range("A23").Value = "Admin"
$wsTables.Range("B23").Value = "Request"
$wsTables.Range("C23").Value = "One to One"

# 3) Rename "Net Banking" table heading to "NetBanking" --------------------
$wsTables.Range("N19").Value = "NetBanking"

# 4) Email Id gets a validity check ------------------------------------------
$wsTables.Range("L13").Value = "CHECK(valid)"

# 5) Branch Code / INT / FOREIGN KEY (row 22 of the Account table) -> red --
$wsTables.Range("F22:H22").Font.Color = 255

# 6) Widen columns B & C to fit the new text in the relationship table ----
$wsTables.Columns.Item(2).ColumnWidth = 11.33203125
$wsTables.Columns.Item(3).ColumnWidth = 12.6640625

# 7) Move the selection, like in the saved file -----------------------------
$wsTables.Activate()
$wsTables.Range("A27").Select()
